$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to
# Text format first, otherwise Excel will coerce them into numeric values
# (losing formatting such as trailing zeros, e.g. "599.26" vs "153.00").
$textCells = @("D5", "D8", "D12", "D16", "D20", "D21", "D23", "D26", "D27", "D31", "D36", "D38", "D39", "D40", "D41", "D42", "D43", "D46", "D48", "D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values from the crypto-price refresh.
$ws.Range("D2").Value = "67.505.48"
$ws.Range("E2").Value = "  -1.46%  "
$ws.Range("D3").Value = "2.676.13"
$ws.Range("E3").Value = "  -0.60%  "
$ws.Range("D5").Value = "599.26"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("E6").Value = "  +3.37%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "0.546"
$ws.Range("E8").Value = "  +0.41%  "
$ws.Range("D9").Value = "2.674.34"
$ws.Range("E9").Value = "  -0.64%  "
$ws.Range("E10").Value = "  +0.19%  "
$ws.Range("E11").Value = "  +1.07%  "
$ws.Range("D12").Value = "0.358"
$ws.Range("E12").Value = "  -0.82%  "
$ws.Range("E13").Value = "  -1.44%  "
$ws.Range("E14").Value = "  -1.64%  "
$ws.Range("D15").Value = "3.162.46"
$ws.Range("E15").Value = "  -0.65%  "
$ws.Range("D16").Value = "0.0000184"
$ws.Range("E16").Value = "  -2.47%  "
$ws.Range("D17").Value = "67.411.50"
$ws.Range("E17").Value = "  -1.57%  "
$ws.Range("D18").Value = "2.699.77"
$ws.Range("E18").Value = "  -0.02%  "
$ws.Range("E19").Value = "  -1.17%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "7.61"
$ws.Range("E20").Value = "  -0.37%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "363.37"
$ws.Range("E21").Value = "  -0.82%  "
$ws.Range("E22").Value = "  -3.64%  "
$ws.Range("D23").Value = "4.81"
$ws.Range("E23").Value = "  -1.45%  "
$ws.Range("E24").Value = "  -3.72%  "
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("D26").Value = "71.20"
$ws.Range("E26").Value = "  -4.45%  "
$ws.Range("D27").Value = "10.23"
$ws.Range("E27").Value = "  +2.02%  "
$ws.Range("D28").Value = "2.806.05"
$ws.Range("E28").Value = "  -0.84%  "
$ws.Range("E29").Value = "  -2.11%  "
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("D31").Value = "554.92"
$ws.Range("E31").Value = "  -2.92%  "
$ws.Range("E32").Value = "  -2.67%  "
$ws.Range("E33").Value = "  -4.24%  "
$ws.Range("E34").Value = "  -0.74%  "
$ws.Range("E35").Value = "  -1.91%  "
$ws.Range("D36").Value = "0.998"
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("E37").Value = "  -4.96%  "
$ws.Range("D38").Value = "19.53"
$ws.Range("E38").Value = "  -1.72%  "
$ws.Range("D39").Value = "154.45"
$ws.Range("E39").Value = "  -4.53%  "
$ws.Range("D40").Value = "0.374"
$ws.Range("E40").Value = "  -1.24%  "
$ws.Range("D41").Value = "5.31"
$ws.Range("E41").Value = "  -1.66%  "
$ws.Range("D42").Value = "1.82"
$ws.Range("E42").Value = "  -4.45%  "
$ws.Range("D43").Value = "17.94"
$ws.Range("E43").Value = "  +0.35%  "
$ws.Range("E44").Value = "  -4.94%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").Value = "40.37"
$ws.Range("E46").Value = "  -0.53%  "
$ws.Range("D47").Value = "0.0₆0297"
$ws.Range("E47").Value = "  -5.95%  "
$ws.Range("D48").Value = "0.592"
$ws.Range("E48").Value = "  -1.10%  "
$ws.Range("D49").Value = "153.00"
$ws.Range("E49").Value = "  -2.79%  "
$ws.Range("E50").Value = "  -3.95%  "
$ws.Range("E51").Value = "  -2.80%  "
